$wb = $excel.ActiveWorkbook

# Add a new worksheet "Sheet2" right after the existing Sheet1
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Populate the data for Sheet2
$ws2.Range("A1").Value = 41234
$ws2.Range("B1").Value = 2131321
$ws2.Range("A2").Value = 1321
$ws2.Range("B2").Value = 321312

# Selection on the new sheet
$ws2.Range("B1").Select()
